$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column). This shifts the
# existing N/O/P columns (Late, Outstanding, Disbursement) one column to the
# right, becoming O/P/Q, and leaves the new N column empty.
$ws.Columns("N").Insert()

# Match the inserted column's width to the width used elsewhere in the sheet.
$ws.Columns("N").ColumnWidth = 9.83

# Update the active selection to reflect where the user ended up after the edit.
$ws.Range("S7").Select()
